$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.617.92"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "'3.318.72"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'580.80"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'173.88"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").Value = "'3.314.90"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").Value = "'0.180"
$ws.Range("E10").Value = "  +3.41%  "
$ws.Range("D11").Value = "'0.577"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "'46.59"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'694.79"
$ws.Range("E14").Value = "  +4.25%  "
$ws.Range("D15").Value = "'3.860.74"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").Value = "'8.34"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "'67.636.21"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'3.325.00"
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("D20").Value = "'17.43"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'11.08"
$ws.Range("E21").Value = "  +2.90%  "
$ws.Range("D22").Value = "'0.887"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "'5.41"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "'16.80"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").Value = "'101.27"
$ws.Range("E25").Value = "  +3.98%  "
$ws.Range("D26").Value = "'3.88"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").Value = "'9.33"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").Value = "'32.83"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "'8.47"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("D31").Value = "'6.94"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "'568.88"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'10.96"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'57.21"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").Value = "'3.704.92"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").Value = "'3.24"
$ws.Range("E38").Value = "  -4.04%  "
$ws.Range("D39").Value = "'34.82"
$ws.Range("E39").Value = "  +7.67%  "
$ws.Range("D40").Value = "'0.132"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "'3.12"
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("D42").Value = "'2.59"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("E43").Value = "  +3.77%  "
$ws.Range("D44").Value = "'0.332"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("D45").Value = "'0.0₃0662"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("E47").Value = "  +1.67%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").Value = "'131.05"
$ws.Range("E51").Value = "  +1.63%  "
